$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Workbook window size (cosmetic view setting captured in the diff)
# ---------------------------------------------------------------------------
$excel.Width = 18360
$excel.Height = 5550

# ---------------------------------------------------------------------------
# 2) Column A: STT now holds plain sequence numbers (stored as text) instead
#    of the SignIn_TCn labels.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "'1"
$ws.Range("A3").Value = "'2"
$ws.Range("A4").Value = "'3"
$ws.Range("A5").Value = "'4"
$ws.Range("A6").Value = "'5"
$ws.Range("A7").Value = "'6"

# Re-apply the original column-A look (border + wrap + vertical-center) that
# plain value assignment can disturb.
$ws.Range("C3").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Column E (Results): every test row is now marked "Pass".
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "Pass"
$ws.Range("E3").Value = "Pass"
$ws.Range("E4").Value = "Pass"
$ws.Range("E5").Value = "Pass"
$ws.Range("E6").Value = "Pass"
$ws.Range("E7").Value = "Pass"

# ---------------------------------------------------------------------------
# 4) Column F: header + data repurposed from "Browser"/"Mozilla" to the new
#    "Expectation key" expected-xpath column.
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Expectation key"
$ws.Range("F2").Value = " //span[text()='Install App']"
$ws.Range("F3").Value = "//span[text()='Incorrect username or password.']"
$ws.Range("F4").Value = "//span[text()='Incorrect username or password.']"
$ws.Range("F5").Value = "//span[text()='Incorrect username or password.']"
$ws.Range("F6").Value = "//span[text()='Incorrect username or password.']"
$ws.Range("F7").Value = "//span[text()='Incorrect username or password.']"

# ---------------------------------------------------------------------------
# 5) Column G: brand new "Expectation result" column (header only, the data
#    cells stay empty but keep a bordered / wrap-text look).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = "Expectation result"

$ws.Range("C3").Copy()
$ws.Range("G2:G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen F & G to fit the new long expectation strings.
$ws.Columns.Item(6).ColumnWidth = 34.67
$ws.Columns.Item(7).ColumnWidth = 37.83

# ---------------------------------------------------------------------------
# 6) Row heights: rows 3-7 grow to fit the two-line descriptions (row 2 stays
#    single-line).
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30

# ---------------------------------------------------------------------------
# 7) Selection cursor ends up on E10 in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("E10").Select()
